$d = $word.ActiveDocument

$d.Content.Find.Execute("Ref-XY7Z9A", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-u248305", 2)

$d.Content.Find.Execute("Ref-DJ49F2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s405246", 2)

$d.Content.Find.Execute("(Ref-A1B2C3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Nguyen, 2015)", 2)

$d.Content.Find.Execute("Ref-AB12CD", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s071588", 2)
